$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 300
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -130
$ws.Range("N12").ClearContents()
$ws.Range("H20").Value = 9999
$ws.Range("I20").Value = 9999
$ws.Range("K20").Value = 9999
$ws.Range("M20").Value = -9769
$ws.Range("H29").Value = 127.25
$ws.Range("I29").Value = 3
$ws.Range("J29").Value = 500
$ws.Range("K29").Value = 9
$ws.Range("L29").Value = 1500
$ws.Range("M29").Value = 272
$ws.Range("N29").Value = -2062
$ws.Range("H31").Value = 33333394
$ws.Range("I31").Value = 37037084
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 111111252
$ws.Range("L31").Value = 600
$ws.Range("M31").Value = -111111022
$ws.Range("N31").Value = -1060
$ws.Range("H34").Value = 5824.6665
$ws.Range("I34").Value = 5824.6665
$ws.Range("K34").Value = 5824.6665
$ws.Range("M34").Value = -5621.6665
$ws.Range("H35").Value = 9999
$ws.Range("I35").Value = 9999
$ws.Range("K35").Value = 9999
$ws.Range("M35").Value = -9620
$ws.Range("H36").Value = 5824.6665
$ws.Range("I36").Value = 5824.6665
$ws.Range("K36").Value = 5824.6665
$ws.Range("M36").Value = -5109.6665
$ws.Range("H41").Value = 382.44446
$ws.Range("I41").Value = 391.7143
$ws.Range("J41").Value = 350
$ws.Range("K41").Value = 391.7143
$ws.Range("L41").Value = 350
$ws.Range("M41").Value = 48.28570000000002
$ws.Range("N41").Value = -1230
$ws.Range("H53").Value = 215.19048
$ws.Range("I53").Value = 150.15384
$ws.Range("K53").Value = 150.15384
$ws.Range("M53").Value = 486.84616
$ws.Range("H62").Value = 7918.3335
$ws.Range("I62").Value = 7873
$ws.Range("J62").Value = 7927.4
$ws.Range("K62").Value = 7873
$ws.Range("L62").Value = 7927.4
$ws.Range("M62").Value = -7249
$ws.Range("N62").Value = -9175.4
$ws.Range("H65").Value = 7918.3335
$ws.Range("I65").Value = 7873
$ws.Range("J65").Value = 7927.4
$ws.Range("K65").Value = 39365
$ws.Range("L65").Value = 39637
$ws.Range("M65").Value = -36245
$ws.Range("N65").Value = -45877
$ws.Range("H107").Value = 802.5
$ws.Range("I107").Value = 802.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 802.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1117.5
$ws.Range("N107").ClearContents()
$ws.Range("H137").Value = 1838.3636
$ws.Range("J137").Value = 1500
$ws.Range("L137").Value = 4500
$ws.Range("N137").Value = -9600

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12532.546
$ws.Range("I32").Value = 13096.556
$ws.Range("K32").Value = 13096.556
$ws.Range("M32").Value = -12809.556
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("H101").Value = 350301
$ws.Range("J101").Value = 350301
$ws.Range("L101").Value = 350301
$ws.Range("N101").Value = -356791
$ws.Range("H102").Value = 999
$ws.Range("I102").Value = 999
$ws.Range("K102").Value = 999
$ws.Range("M102").Value = 623
$ws.Range("H119").Value = 51199.6
$ws.Range("J119").Value = 51199.6
$ws.Range("L119").Value = 51199.6
$ws.Range("N119").Value = -60875.6
$ws.Range("H132").Value = 3999.6667
$ws.Range("I132").Value = 3999
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 11997
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -9467
$ws.Range("N132").Value = -17060
$ws.Range("H135").Value = 197499.33
$ws.Range("J135").Value = 197499.33
$ws.Range("L135").Value = 197499.33
$ws.Range("N135").Value = -207639.33

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1967.8572
$ws.Range("I94").Value = 1965.8334
$ws.Range("K94").Value = 1965.8334
$ws.Range("M94").Value = -1514.8334
$ws.Range("H99").Value = 2091.4546
$ws.Range("I99").Value = 2111.6
$ws.Range("K99").Value = 2111.6
$ws.Range("M99").Value = -613.5999999999999
$ws.Range("H105").Value = 3349
$ws.Range("I105").Value = 3999
$ws.Range("K105").Value = 3999
$ws.Range("M105").Value = -2252
$ws.Range("H119").Value = 39994
$ws.Range("J119").Value = 39994
$ws.Range("L119").Value = 39994
$ws.Range("N119").Value = -49670

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 128.33333
$ws.Range("J7").Value = 193.63637
$ws.Range("L7").Value = 193.63637
$ws.Range("N7").Value = -419.63637
$ws.Range("H26").Value = 3339.6667
$ws.Range("I26").Value = 1009.5
$ws.Range("J26").Value = 8000
$ws.Range("K26").Value = 1009.5
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = -722.5
$ws.Range("N26").Value = -8574

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 79984.336
$ws.Range("J37").Value = 79984.336
$ws.Range("L37").Value = 239953.008
$ws.Range("N37").Value = -240177.008

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 113.4
$ws.Range("I2").Value = 135.75
$ws.Range("J2").Value = 24
$ws.Range("K2").Value = 135.75
$ws.Range("L2").Value = 24
$ws.Range("M2").Value = -22.75
$ws.Range("N2").Value = -250
$ws.Range("H70").Value = 7981.5884
$ws.Range("J70").Value = 8112.857
$ws.Range("L70").Value = 8112.857
$ws.Range("N70").Value = -8652.857
$ws.Range("H73").Value = 7981.5884
$ws.Range("J73").Value = 8112.857
$ws.Range("L73").Value = 8112.857
$ws.Range("N73").Value = -9984.857
$ws.Range("H102").Value = 1084.2727
$ws.Range("I102").Value = 1084.2727
$ws.Range("K102").Value = 1084.2727
$ws.Range("M102").Value = 537.7273
$ws.Range("H122").Value = 4214
$ws.Range("I122").Value = 2073.6667
$ws.Range("K122").Value = 6221.000100000001
$ws.Range("M122").Value = -3771.000100000001
$ws.Range("H126").Value = 3940
$ws.Range("I126").Value = 3880
$ws.Range("K126").Value = 11640
$ws.Range("M126").Value = -9170
$ws.Range("H132").Value = 966.3333
$ws.Range("I132").Value = 950
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 2850
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -320
$ws.Range("N132").Value = -8057

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1413.1428
$ws.Range("I16").Value = 958.6
$ws.Range("J16").Value = 2549.5
$ws.Range("K16").Value = 958.6
$ws.Range("L16").Value = 2549.5
$ws.Range("M16").Value = -788.6
$ws.Range("N16").Value = -2889.5
$ws.Range("H22").Value = 1639.4
$ws.Range("J22").Value = 1799.75
$ws.Range("L22").Value = 1799.75
$ws.Range("N22").Value = -2389.75
$ws.Range("H27").Value = 1639.4
$ws.Range("J27").Value = 1799.75
$ws.Range("L27").Value = 1799.75
$ws.Range("N27").Value = -2013.75
$ws.Range("H40").Value = 4054.111
$ws.Range("I40").Value = 2300
$ws.Range("K40").Value = 2300
$ws.Range("M40").Value = -2164
$ws.Range("H46").Value = 33760.75
$ws.Range("I46").Value = 64521.625
$ws.Range("K46").Value = 64521.625
$ws.Range("M46").Value = -64333.625
$ws.Range("H55").Value = 241.85
$ws.Range("I55").Value = 185.38889
$ws.Range("J55").Value = 750
$ws.Range("K55").Value = 185.38889
$ws.Range("L55").Value = 750
$ws.Range("M55").Value = -12.38889
$ws.Range("N55").Value = -1096
$ws.Range("H61").Value = 4134.4287
$ws.Range("I61").Value = 4157
$ws.Range("J61").Value = 3999
$ws.Range("K61").Value = 4157
$ws.Range("L61").Value = 3999
$ws.Range("M61").Value = -3955
$ws.Range("N61").Value = -4403
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H113").Value = 4134.4287
$ws.Range("I113").Value = 4157
$ws.Range("J113").Value = 3999
$ws.Range("K113").Value = 4157
$ws.Range("L113").Value = 3999
$ws.Range("M113").Value = -1987
$ws.Range("N113").Value = -8339

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1999
$ws.Range("I107").Value = 1999
$ws.Range("K107").Value = 5997
$ws.Range("M107").Value = -4077
$ws.Range("H119").Value = 99999.664
$ws.Range("J119").Value = 99999.664
$ws.Range("L119").Value = 99999.664
$ws.Range("N119").Value = -109675.664
$ws.Range("H122").Value = 14167.667
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 20751.5
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 62254.5
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -67154.5
$ws.Range("H126").Value = 4539.85
$ws.Range("I126").Value = 3358.5833
$ws.Range("J126").Value = 6311.75
$ws.Range("K126").Value = 10075.7499
$ws.Range("L126").Value = 18935.25
$ws.Range("M126").Value = -7605.749899999999
$ws.Range("N126").Value = -23875.25
$ws.Range("H136").Value = 2152.6
$ws.Range("I136").Value = 1949.2142
$ws.Range("K136").Value = 5847.642599999999
$ws.Range("M136").Value = -3297.642599999999

